$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.329.06'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '2.375.48'
$ws.Range('E3').Value = '  +3.08%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''309.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').Value = '''104.61'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.62%  '
$ws.Range('E7').Value = '  -1.71%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '''0.519'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.24%  '
$ws.Range('D10').Value = '''36.30'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.89%  '
$ws.Range('D11').Value = '''52.92'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.73%  '
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('E13').Value = '  -0.71%  '
$ws.Range('E14').Value = '  +0.70%  '
$ws.Range('D15').Value = '2.740.27'
$ws.Range('E15').Value = '  +2.98%  '
$ws.Range('D16').Value = '''15.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.29%  '
$ws.Range('D17').Value = '2.374.63'
$ws.Range('E17').Value = '  +3.45%  '
$ws.Range('E18').Value = '  +1.84%  '
$ws.Range('D19').Value = '43.308.95'
$ws.Range('E19').Value = '  +0.60%  '
$ws.Range('E20').Value = '  -4.16%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0927'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '''6.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.82%  '
$ws.Range('D23').Value = '''68.48'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('D24').Value = '''242.18'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.96%  '
$ws.Range('E25').Value = '  +1.90%  '
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('D28').Value = '''26.11'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.30%  '
$ws.Range('E29').Value = '  +8.72%  '
$ws.Range('D30').Value = '''36.90'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.01%  '
$ws.Range('D31').Value = '''9.62'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('D32').Value = '''161.87'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.13%  '
$ws.Range('D33').Value = '''5.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.65%  '
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('D35').Value = '''18.39'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.34%  '
$ws.Range('E36').Value = '  +6.58%  '
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').Value = '''1.96'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.98%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''4.68'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +11.80%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '''0.0741'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('D43').Value = '''2.44'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.37%  '
$ws.Range('D44').Value = '''20.30'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.67%  '
$ws.Range('D45').Value = '2.005.71'
$ws.Range('E45').Value = '  +1.84%  '
$ws.Range('D46').Value = '''3.22'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.61%  '
$ws.Range('D47').Value = '''0.0292'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('E48').Value = '  +5.90%  '
$ws.Range('D49').Value = '''58.01'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.35%  '
$ws.Range('D50').Value = '''2.95'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.72%  '
$ws.Range('D51').Value = '2.578.26'
$ws.Range('E51').Value = '  +1.90%  '
